$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header rich-text strings (Volume number, report week dates) ---
$a8 = $ws.Range("A8")
$a8full = $a8.Value2
$a8idx = $a8full.LastIndexOf("19") + 1
$a8.Characters($a8idx, 2).Text = "20"

$c9 = $ws.Range("C9")
$c9full = $c9.Value2
$c9idx1 = $c9full.IndexOf("5/8/2023") + 1
$c9.Characters($c9idx1, 8).Text = "5/15/2023"
$c9full2 = $ws.Range("C9").Value2
$c9idx2 = $c9full2.IndexOf("5/14/2023") + 1
$ws.Range("C9").Characters($c9idx2, 9).Text = "5/21/2023"

# --- Update crime-stat data table (rows 15-27) ---
$ws.Range("J14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("N14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("J14").Copy($ws.Range("G15"))
$ws.Range("G15").Value = 1
$ws.Range("N14").Copy($ws.Range("H15"))
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 20
$ws.Range("N15").Value = -50
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -40
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 50
$ws.Range("M16").Value = -25
$ws.Range("N16").Value = -86.899563318777
$ws.Range("C17").Value = 3
$ws.Range("C14").Copy($ws.Range("D17"))
$ws.Range("E14").Copy($ws.Range("E17"))
$ws.Range("I17").Value = 58
$ws.Range("K17").Value = 45
$ws.Range("L17").Value = 75.757575757575
$ws.Range("M17").Value = 52.631578947368
$ws.Range("N17").Value = -39.583333333333
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 46
$ws.Range("K18").Value = 64.285714285714
$ws.Range("L18").Value = 31.428571428571
$ws.Range("M18").Value = -52.083333333333
$ws.Range("N18").Value = -89.61625282167
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 8.108108108108
$ws.Range("I19").Value = 163
$ws.Range("J19").Value = 205
$ws.Range("K19").Value = -20.487804878048
$ws.Range("L19").Value = 29.365079365079
$ws.Range("M19").Value = 44.247787610619
$ws.Range("N19").Value = -11.41304347826
$ws.Range("J14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 46
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = 12.195121951219
$ws.Range("L20").Value = 130
$ws.Range("M20").Value = -24.590163934426
$ws.Range("N20").Value = -94.162436548223
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -5.263157894736
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = 1.492537313432
$ws.Range("I21").Value = 351
$ws.Range("J21").Value = 341
$ws.Range("K21").Value = 2.932551319648
$ws.Range("L21").Value = 48.101265822784
$ws.Range("M21").Value = -0.847457627118
$ws.Range("N21").Value = -79.988597491448
$ws.Range("J14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("N14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -57.142857142857
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 27.272727272727
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -5.263157894736
$ws.Range("I24").Value = 608
$ws.Range("J24").Value = 704
$ws.Range("K24").Value = -13.636363636363
$ws.Range("L24").Value = 94.871794871794
$ws.Range("M24").Value = 36.322869955157
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 25
$ws.Range("I25").Value = 145
$ws.Range("J25").Value = 129
$ws.Range("K25").Value = 12.403100775193
$ws.Range("L25").Value = 57.608695652173
$ws.Range("M25").Value = 12.403100775193
$ws.Range("J14").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("N14").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("J14").Copy($ws.Range("G26"))
$ws.Range("G26").Value = 1
$ws.Range("N14").Copy($ws.Range("H26"))
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = 33.333333333333
$ws.Range("L26").Value = 60
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
